# Generate Report for Handback
# The handback for "ecf43ef3-1d62-4099-80e8-e264b6095bb6.md" has completed: update the
# status rows across the Overview / zh-cn / de-de sheets to reflect a successful handback
# instead of the previous "Ready for handoff" / stale-handback-error state.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the ecf43ef3 file ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the ecf43ef3 file ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-22 04:57:31"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 13.7

# --- de-de sheet: row 3 is the ecf43ef3 file ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-22 04:57:37"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 13.7
